$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.935.94'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '2.550.12'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.98'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.30%  '
$ws.Range('E7').Value = '  +0.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.547'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.94'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0833'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.12%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.76'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.52%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.115'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('D14').Value = '2.934.75'
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('D15').Value = '2.556.27'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.11'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.875'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '42.906.20'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.06%  '
$ws.Range('D20').Value = '0.0₃0995'
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.09'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '257.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.09'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '28.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.01%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.43%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.78%  '
$ws.Range('E31').Value = '  +3.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '158.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.56'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +15.50%  '
$ws.Range('E34').Value = '  -2.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0802'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('E37').Value = '  -4.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.117'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.02'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.53%  '
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.45'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.92'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.43%  '
$ws.Range('E43').Value = '  +29.02%  '
$ws.Range('D44').Value = '2.102.04'
$ws.Range('E44').Value = '  +0.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0306'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.80%  '
$ws.Range('E48').Value = '  -1.52%  '
$ws.Range('D49').Value = '2.792.37'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.03'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.05%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.41%  '
